# Change model to closely match handle-API
#
# - Rename/repurpose "PID4CatRecord" -> "HandleAPIRecord" with new columns
# - Insert three new sheets: HandleRecord, HandleData, HandleRecordContainer
#   right after HandleAPIRecord
# - Remove the old "Container" sheet (its single column now lives in
#   the new "HandleRecordContainer" sheet)
# - All other sheets (PID4CatRelation, ResourceInfo, LogRecord, Agent,
#   RepresentationVariant) are left untouched.
#
# NOTE: new sheet references are always re-fetched by name via
# $wb.Worksheets.Item("...") instead of cached in variables, since the
# After/Before positional lookup of a stored reference to a just-renamed
# worksheet does not behave reliably in this runtime.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- 1. Repurpose the first sheet: PID4CatRecord -> HandleAPIRecord -------
# drop the old dropdown validation on column C (status) before reshaping
$wb.Worksheets.Item("PID4CatRecord").Range("C2:C1048576").Validation.Delete()

# clear the old header row (A1:I1) then write the new, smaller header
$wb.Worksheets.Item("PID4CatRecord").Range("A1:I1").ClearContents()
$wb.Worksheets.Item("PID4CatRecord").Cells.Item(1, 1).Value = "response_code"
$wb.Worksheets.Item("PID4CatRecord").Cells.Item(1, 2).Value = "handle"
$wb.Worksheets.Item("PID4CatRecord").Cells.Item(1, 3).Value = "values"

$wb.Worksheets.Item("PID4CatRecord").Name = "HandleAPIRecord"

# Use the existing single-column "Container" sheet as a style/structure
# template for the new sheets, so they inherit the same sheetPr/margins as
# the rest of the workbook; then rename + re-populate each copy.

# --- 2. New sheet: HandleRecord --------------------------------------------
$wb.Worksheets.Item("Container").Copy($null, $wb.Worksheets.Item("Container"))
$wb.Worksheets.Item("Container (2)").Name = "HandleRecord"

$wb.Worksheets.Item("HandleRecord").Range("A1:A1").ClearContents()
$wb.Worksheets.Item("HandleRecord").Cells.Item(1, 1).Value = "index"
$wb.Worksheets.Item("HandleRecord").Cells.Item(1, 2).Value = "type"
$wb.Worksheets.Item("HandleRecord").Cells.Item(1, 3).Value = "data"
$wb.Worksheets.Item("HandleRecord").Cells.Item(1, 4).Value = "ttl"
$wb.Worksheets.Item("HandleRecord").Cells.Item(1, 5).Value = "timestamp"
$wb.Worksheets.Item("HandleRecord").Range("B2:B1048576").Validation.Add(3, 1, 1, """URL,STATUS,SCHEMA_VER,LICENSE,EMAIL,RESOURCE_INFO,RELATED,LOG""")

# --- 3. New sheet: HandleData ----------------------------------------------
$wb.Worksheets.Item("Container").Copy($null, $wb.Worksheets.Item("Container"))
$wb.Worksheets.Item("Container (2)").Name = "HandleData"

$wb.Worksheets.Item("HandleData").Range("A1:A1").ClearContents()
$wb.Worksheets.Item("HandleData").Cells.Item(1, 1).Value = "format"
$wb.Worksheets.Item("HandleData").Cells.Item(1, 2).Value = "value"

# --- 4. New sheet: HandleRecordContainer -----------------------------------
$wb.Worksheets.Item("Container").Copy($null, $wb.Worksheets.Item("Container"))
$wb.Worksheets.Item("Container (2)").Name = "HandleRecordContainer"
# content ("contains_pids") is already correct, inherited from the template

# --- 5. Reorder the three new sheets to sit right after HandleAPIRecord ---
$wb.Worksheets.Item("HandleRecordContainer").Move($wb.Worksheets.Item("PID4CatRelation"))
$wb.Worksheets.Item("HandleData").Move($wb.Worksheets.Item("HandleRecordContainer"))
$wb.Worksheets.Item("HandleRecord").Move($wb.Worksheets.Item("HandleData"))

# --- 6. Remove the old Container sheet (superseded by HandleRecordContainer)
$wb.Worksheets.Item("Container").Delete() | Out-Null

# keep the first sheet as the active one, as in the original workbook
$wb.Worksheets.Item(1).Activate()
